# Weekly update: insert the new "Ajo" (garlic) price record for
# Feria Lagunitas de Puerto Montt at the top of the date block (row 348),
# pushing the existing rows 348:383 down to 349:384.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 348; everything at/after 348 shifts down one row.
$ws.Rows.Item(348).Insert()

# Populate the new row with the latest weekly observation.
$ws.Cells.Item(348, 1).Value  = 4
$ws.Cells.Item(348, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(348, 3).Value  = "Los Lagos"
$ws.Cells.Item(348, 4).Value  = 44918
$ws.Cells.Item(348, 5).Value  = 10
$ws.Cells.Item(348, 6).Value  = 100112003
$ws.Cells.Item(348, 7).Value  = "Ajo"
$ws.Cells.Item(348, 8).Value  = "Chino"
$ws.Cells.Item(348, 9).Value  = "Primera"
$ws.Cells.Item(348, 10).Value = 80
$ws.Cells.Item(348, 11).Value = 18000
$ws.Cells.Item(348, 12).Value = 18000
$ws.Cells.Item(348, 13).Value = 18000
$ws.Cells.Item(348, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(348, 15).Value = "China"
$ws.Cells.Item(348, 16).Value = 1800
$ws.Cells.Item(348, 17).Value = 10
$ws.Cells.Item(348, 18).Value = "Hortaliza"
